# Update cryptos list figures (prices + 1h volume change %) as scraped on
# Wed Sep 11 19:28:09 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (new Price text, new Volume(1h) text). A $null entry means
# that column is left untouched for that row. Price values that would be
# auto-recognized by Excel as a plain number (single decimal point) are
# prefixed with a leading apostrophe so they are stored/forced as text,
# exactly like the source data (which keeps these as text cells).
$updates = @(
    @{ Row = 2;  D = "57.075.77";  E = "  -0.40%  " },
    @{ Row = 3;  D = "2.315.59";   E = "  -1.44%  " },
    @{ Row = 4;  D = $null;        E = "  -0.01%  " },
    @{ Row = 5;  D = "'532.37";    E = "  +2.31%  " },
    @{ Row = 6;  D = "'131.86";    E = "  -3.07%  " },
    @{ Row = 7;  D = $null;        E = "  -0.31%  " },
    @{ Row = 8;  D = $null;        E = "  -0.24%  " },
    @{ Row = 9;  D = "2.336.90";   E = "  -1.08%  " },
    @{ Row = 10; D = $null;        E = "  -1.13%  " },
    @{ Row = 11; D = $null;        E = "  +0.23%  " },
    @{ Row = 12; D = $null;        E = "  -2.50%  " },
    @{ Row = 13; D = $null;        E = "  +0.49%  " },
    @{ Row = 14; D = "2.733.99";   E = "  -1.19%  " },
    @{ Row = 15; D = "'23.43";     E = "  -3.23%  " },
    @{ Row = 16; D = "57.089.05";  E = "  -0.34%  " },
    @{ Row = 17; D = $null;        E = "  -2.16%  " },
    @{ Row = 18; D = "2.338.21";   E = "  -0.12%  " },
    @{ Row = 19; D = "'337.51";    E = "  +2.69%  " },
    @{ Row = 20; D = "'10.44";     E = "  -1.50%  " },
    @{ Row = 21; D = "'6.90";      E = "  +2.78%  " },
    @{ Row = 22; D = $null;        E = "  -2.12%  " },
    @{ Row = 23; D = $null;        E = "  +0.07%  " },
    @{ Row = 24; D = "'61.61";     E = "  +0.47%  " },
    @{ Row = 25; D = "'8.74";      E = "  +5.26%  " },
    @{ Row = 26; D = $null;        E = "  +0.67%  " },
    @{ Row = 27; D = $null;        E = "  -0.75%  " },
    @{ Row = 28; D = $null;        E = "  +0.77%  " },
    @{ Row = 29; D = $null;        E = "  +0.17%  " },
    @{ Row = 30; D = "'1.73";      E = "  +1.62%  " },
    @{ Row = 31; D = $null;        E = "  -3.20%  " },
    @{ Row = 32; D = $null;        E = "  -2.88%  " },
    @{ Row = 33; D = $null;        E = "  -0.31%  " },
    @{ Row = 35; D = $null;        E = "  -0.26%  " },
    @{ Row = 36; D = $null;        E = "  -2.99%  " },
    @{ Row = 37; D = "'3.99";      E = "  -1.22%  " },
    @{ Row = 38; D = "'0.903";     E = "  -2.32%  " },
    @{ Row = 39; D = $null;        E = "  +0.35%  " },
    @{ Row = 40; D = $null;        E = "  +1.53%  " },
    @{ Row = 41; D = $null;        E = "  +9.68%  " },
    @{ Row = 42; D = "'148.68";    E = "  -1.25%  " },
    @{ Row = 43; D = $null;        E = "  -1.41%  " },
    @{ Row = 44; D = $null;        E = "  -1.50%  " },
    @{ Row = 45; D = "'279.89";    E = "  -0.89%  " },
    @{ Row = 46; D = "'0.0927";    E = $null },
    @{ Row = 47; D = $null;        E = "  -1.41%  " },
    @{ Row = 48; D = $null;        E = "  +2.12%  " },
    @{ Row = 49; D = $null;        E = "  -1.37%  " },
    @{ Row = 50; D = $null;        E = "  +0.23%  " },
    @{ Row = 51; D = $null;        E = "  -2.25%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $ws.Cells.Item($r, 4).Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}
